# Updated cryptos list (GitHub Actions scheduled refresh): new Price (col D)
# and Volume(1h) (col E) text values for rows 2-51. Price cells that look
# like plain numbers are pre-formatted as Text ("@") before assignment so
# Excel keeps them as literal strings (e.g. "45.53") instead of silently
# coercing them to numeric cells, matching the source data's inline-string
# representation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.563.20'
$ws.Range('E2').Value = '  -2.51%  '
$ws.Range('D3').Value = '2.223.89'
$ws.Range('E3').Value = '  -2.25%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '112.94'
$ws.Range('E5').Value = '  -7.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '297.65'
$ws.Range('E6').Value = '  +12.23%  '
$ws.Range('E7').Value = '  -1.43%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.615'
$ws.Range('E9').Value = '  -0.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '45.53'
$ws.Range('E10').Value = '  -5.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0928'
$ws.Range('E11').Value = '  -1.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.80'
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.94'
$ws.Range('E13').Value = '  -0.59%  '
$ws.Range('E14').Value = '  -2.96%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.944'
$ws.Range('E15').Value = '  +5.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.19'
$ws.Range('E16').Value = '  -1.58%  '
$ws.Range('D17').Value = '2.559.57'
$ws.Range('E17').Value = '  -2.35%  '
$ws.Range('D18').Value = '2.264.51'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').Value = '42.441.27'
$ws.Range('E19').Value = '  -2.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.41'
$ws.Range('E20').Value = '  +6.42%  '
$ws.Range('E21').Value = '  -2.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.73'
$ws.Range('E22').Value = '  +2.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.58'
$ws.Range('E23').Value = '  +25.25%  '
$ws.Range('E24').Value = '  -4.77%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '230.51'
$ws.Range('E25').Value = '  -1.95%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.51'
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.85'
$ws.Range('E27').Value = '  +0.64%  '
$ws.Range('E28').Value = '  -1.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.91'
$ws.Range('E29').Value = '  -1.45%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '39.33'
$ws.Range('E30').Value = '  -6.26%  '
$ws.Range('E31').Value = '  -1.59%  '
$ws.Range('E32').Value = '  -3.57%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '174.08'
$ws.Range('E33').Value = '  +1.46%  '
$ws.Range('E34').Value = '  -1.87%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0893'
$ws.Range('E35').Value = '  -1.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.80'
$ws.Range('E36').Value = '  +1.42%  '
$ws.Range('E37').Value = '  +6.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.89'
$ws.Range('E38').Value = '  +4.69%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.127'
$ws.Range('E39').Value = '  -1.48%  '
$ws.Range('E40').Value = '  -2.17%  '
$ws.Range('E41').Value = '  -2.96%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.49'
$ws.Range('E42').Value = '  -1.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.241'
$ws.Range('E43').Value = '  +1.53%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '71.48'
$ws.Range('E44').Value = '  -4.72%  '
$ws.Range('E45').Value = '  -4.57%  '
$ws.Range('E46').Value = '  +0.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.33'
$ws.Range('E47').Value = '  -2.83%  '
$ws.Range('E48').Value = '  -4.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.32'
$ws.Range('E49').Value = '  +4.90%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '105.27'
$ws.Range('E50').Value = '  +3.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.58'
$ws.Range('E51').Value = '  +0.35%  '
